# Generate Report for Handback
# Updates the Latest HO Xliff Generate Date / Correspond Handoff & Handback
# datetimes for 93463dfa-68cf-422f-b64f-7fbecba66435.md after a fresh
# handback report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 4 is 93463dfa-68cf-422f-b64f-7fbecba66435.md
# Column G = "Latest HO Xliff Generate Date"
$overview.Range("G4").Value = "2017-02-21 03:15:25"

# zh-cn sheet: row 4 is 93463dfa-68cf-422f-b64f-7fbecba66435.md
# Column H = "Correspond Handoff Datetime", Column L = "Correspond Handback DateTime"
$zhcn.Range("H4").Value = "2017-02-21 03:15:09"
$zhcn.Range("L4").Value = "2017-02-21 03:16:02"

# de-de sheet: row 4 is 93463dfa-68cf-422f-b64f-7fbecba66435.md
# Column H = "Correspond Handoff Datetime", Column L = "Correspond Handback DateTime"
$dede.Range("H4").Value = "2017-02-21 03:15:25"
$dede.Range("L4").Value = "2017-02-21 03:16:24"
